$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 511.45456
$ws.Cells.Item(33, 9).Value = 97.0625
$ws.Cells.Item(33, 10).Value = 1616.5
$ws.Cells.Item(33, 11).Value = 97.0625
$ws.Cells.Item(33, 12).Value = 1616.5
$ws.Cells.Item(33, 13).Value = 131.9375
$ws.Cells.Item(33, 14).Value = -2074.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 675.55554
$ws.Cells.Item(98, 9).Value = 636.86957
$ws.Cells.Item(98, 11).Value = 636.86957
$ws.Cells.Item(98, 13).Value = 861.13043

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1902.963
$ws.Cells.Item(112, 9).Value = 837.5
$ws.Cells.Item(112, 10).Value = 2088.261
$ws.Cells.Item(112, 11).Value = 2512.5
$ws.Cells.Item(112, 12).Value = 6264.782999999999
$ws.Cells.Item(112, 13).Value = -1404.5
$ws.Cells.Item(112, 14).Value = -8480.782999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 675.55554
$ws.Cells.Item(122, 9).Value = 636.86957
$ws.Cells.Item(122, 11).Value = 1910.60871
$ws.Cells.Item(122, 13).Value = 539.39129

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 878.4138
$ws.Cells.Item(129, 10).Value = 1034.7391
$ws.Cells.Item(129, 12).Value = 3104.2173
$ws.Cells.Item(129, 14).Value = -13104.2173

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1563812.6
$ws.Cells.Item(137, 9).Value = 2128622.5
$ws.Cells.Item(137, 10).Value = 2279.4707
$ws.Cells.Item(137, 11).Value = 6385867.5
$ws.Cells.Item(137, 12).Value = 6838.4121
$ws.Cells.Item(137, 13).Value = -6383317.5
$ws.Cells.Item(137, 14).Value = -11938.4121

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1939794.9
$ws.Cells.Item(138, 9).Value = 1200.0358
$ws.Cells.Item(138, 10).Value = 5558505
$ws.Cells.Item(138, 11).Value = 3600.1074
$ws.Cells.Item(138, 12).Value = 16675515
$ws.Cells.Item(138, 13).Value = 1539.8926
$ws.Cells.Item(138, 14).Value = -16685795

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1015.05
$ws.Cells.Item(32, 9).Value = 801.23334
$ws.Cells.Item(32, 10).Value = 2939.4
$ws.Cells.Item(32, 11).Value = 801.23334
$ws.Cells.Item(32, 12).Value = 2939.4
$ws.Cells.Item(32, 13).Value = -514.23334
$ws.Cells.Item(32, 14).Value = -3513.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 20876110
$ws.Cells.Item(61, 9).Value = 23834410
$ws.Cells.Item(61, 10).Value = 168007
$ws.Cells.Item(61, 11).Value = 23834410
$ws.Cells.Item(61, 12).Value = 168007
$ws.Cells.Item(61, 13).Value = -23834198
$ws.Cells.Item(61, 14).Value = -168431

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 11456548
$ws.Cells.Item(74, 9).Value = 20918620
$ws.Cells.Item(74, 10).Value = 102060
$ws.Cells.Item(74, 11).Value = 20918620
$ws.Cells.Item(74, 12).Value = 102060
$ws.Cells.Item(74, 13).Value = -20917746
$ws.Cells.Item(74, 14).Value = -103808

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 11456548
$ws.Cells.Item(77, 9).Value = 20918620
$ws.Cells.Item(77, 10).Value = 102060
$ws.Cells.Item(77, 11).Value = 104593100
$ws.Cells.Item(77, 12).Value = 510300
$ws.Cells.Item(77, 13).Value = -104588732
$ws.Cells.Item(77, 14).Value = -519036

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(131, 8).Value = 58779.316
$ws.Cells.Item(131, 10).Value = 58779.316
$ws.Cells.Item(131, 12).Value = 58779.316
$ws.Cells.Item(131, 14).Value = -68859.31599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 20876110
$ws.Cells.Item(136, 9).Value = 23834410
$ws.Cells.Item(136, 10).Value = 168007
$ws.Cells.Item(136, 11).Value = 71503230
$ws.Cells.Item(136, 12).Value = 504021
$ws.Cells.Item(136, 13).Value = -71500680
$ws.Cells.Item(136, 14).Value = -509121

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 11795.375
$ws.Cells.Item(86, 9).Value = 14792.277
$ws.Cells.Item(86, 11).Value = 14792.277
$ws.Cells.Item(86, 13).Value = -13669.277

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 11795.375
$ws.Cells.Item(89, 9).Value = 14792.277
$ws.Cells.Item(89, 11).Value = 73961.38499999999
$ws.Cells.Item(89, 13).Value = -68345.38499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1527.0182
$ws.Cells.Item(134, 9).Value = 966.0952
$ws.Cells.Item(134, 10).Value = 3339.2307
$ws.Cells.Item(134, 11).Value = 2898.2856
$ws.Cells.Item(134, 12).Value = 10017.6921
$ws.Cells.Item(134, 13).Value = -363.2856000000002
$ws.Cells.Item(134, 14).Value = -15087.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2667.3262
$ws.Cells.Item(31, 9).Value = 1116.75
$ws.Cells.Item(31, 11).Value = 1116.75
$ws.Cells.Item(31, 13).Value = -821.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2667.3262
$ws.Cells.Item(34, 9).Value = 1116.75
$ws.Cells.Item(34, 11).Value = 1116.75
$ws.Cells.Item(34, 13).Value = -914.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 23257462
$ws.Cells.Item(58, 9).Value = 24391680
$ws.Cells.Item(58, 10).Value = 6000.5
$ws.Cells.Item(58, 11).Value = 24391680
$ws.Cells.Item(58, 12).Value = 6000.5
$ws.Cells.Item(58, 13).Value = -24391477
$ws.Cells.Item(58, 14).Value = -6406.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 32253.092
$ws.Cells.Item(132, 9).Value = 2018.88
$ws.Cells.Item(132, 10).Value = 126735
$ws.Cells.Item(132, 11).Value = 6056.64
$ws.Cells.Item(132, 12).Value = 380205
$ws.Cells.Item(132, 13).Value = -3526.64
$ws.Cells.Item(132, 14).Value = -385265

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 44562.04
$ws.Cells.Item(134, 9).Value = 1229.8334
$ws.Cells.Item(134, 10).Value = 155987.72
$ws.Cells.Item(134, 11).Value = 3689.5002
$ws.Cells.Item(134, 12).Value = 467963.16
$ws.Cells.Item(134, 13).Value = -1154.5002
$ws.Cells.Item(134, 14).Value = -473033.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 23257462
$ws.Cells.Item(136, 9).Value = 24391680
$ws.Cells.Item(136, 10).Value = 6000.5
$ws.Cells.Item(136, 11).Value = 73175040
$ws.Cells.Item(136, 12).Value = 18001.5
$ws.Cells.Item(136, 13).Value = -73172490
$ws.Cells.Item(136, 14).Value = -23101.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 541.0769
$ws.Cells.Item(5, 9).Value = 440.1579
$ws.Cells.Item(5, 10).Value = 815
$ws.Cells.Item(5, 11).Value = 1320.4737
$ws.Cells.Item(5, 12).Value = 2445
$ws.Cells.Item(5, 13).Value = -1208.4737
$ws.Cells.Item(5, 14).Value = -2669

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 587.67645
$ws.Cells.Item(122, 9).Value = 278.05
$ws.Cells.Item(122, 10).Value = 1030
$ws.Cells.Item(122, 11).Value = 2502.45
$ws.Cells.Item(122, 12).Value = 9270
$ws.Cells.Item(122, 13).Value = -52.45000000000027
$ws.Cells.Item(122, 14).Value = -14170

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 2978035.2
$ws.Cells.Item(129, 9).Value = 1330
$ws.Cells.Item(129, 10).Value = 4904138.5
$ws.Cells.Item(129, 11).Value = 3990
$ws.Cells.Item(129, 12).Value = 14712415.5
$ws.Cells.Item(129, 13).Value = 1010
$ws.Cells.Item(129, 14).Value = -14722415.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(130, 8).Value = 2671.111
$ws.Cells.Item(130, 9).Value = 1000
$ws.Cells.Item(130, 10).Value = 3148.5715
$ws.Cells.Item(130, 11).Value = 3000
$ws.Cells.Item(130, 12).Value = 9445.7145
$ws.Cells.Item(130, 13).Value = 2020
$ws.Cells.Item(130, 14).Value = -19485.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1081.7018
$ws.Cells.Item(131, 9).Value = 562.7143
$ws.Cells.Item(131, 10).Value = 1154.36
$ws.Cells.Item(131, 11).Value = 1688.1429
$ws.Cells.Item(131, 12).Value = 3463.08
$ws.Cells.Item(131, 13).Value = 3351.8571
$ws.Cells.Item(131, 14).Value = -13543.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 5735.619
$ws.Cells.Item(133, 9).Value = 2720
$ws.Cells.Item(133, 10).Value = 7997.3335
$ws.Cells.Item(133, 11).Value = 8160
$ws.Cells.Item(133, 12).Value = 23992.0005
$ws.Cells.Item(133, 13).Value = -3100
$ws.Cells.Item(133, 14).Value = -34112.00049999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 5116.9287
$ws.Cells.Item(134, 9).Value = 1523.6666
$ws.Cells.Item(134, 10).Value = 7811.875
$ws.Cells.Item(134, 11).Value = 4570.9998
$ws.Cells.Item(134, 12).Value = 23435.625
$ws.Cells.Item(134, 13).Value = 499.0002000000004
$ws.Cells.Item(134, 14).Value = -33575.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 541.0769
$ws.Cells.Item(135, 9).Value = 440.1579
$ws.Cells.Item(135, 10).Value = 815
$ws.Cells.Item(135, 11).Value = 3961.4211
$ws.Cells.Item(135, 12).Value = 7335
$ws.Cells.Item(135, 13).Value = -1426.4211
$ws.Cells.Item(135, 14).Value = -12405

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value = 2775.8823
$ws.Cells.Item(136, 9).Value = 2302.8572
$ws.Cells.Item(136, 11).Value = 6908.571599999999
$ws.Cells.Item(136, 13).Value = -1808.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 36704.65
$ws.Cells.Item(137, 9).Value = 999.6667
$ws.Cells.Item(137, 10).Value = 44355.715
$ws.Cells.Item(137, 11).Value = 2999.0001
$ws.Cells.Item(137, 12).Value = 133067.145
$ws.Cells.Item(137, 13).Value = 2100.9999
$ws.Cells.Item(137, 14).Value = -143267.145

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 24390
$ws.Cells.Item(62, 10).Value = 24390
$ws.Cells.Item(62, 12).Value = 24390
$ws.Cells.Item(62, 14).Value = -25762

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(65, 8).Value = 24390
$ws.Cells.Item(65, 10).Value = 24390
$ws.Cells.Item(65, 12).Value = 73170
$ws.Cells.Item(65, 14).Value = -80034

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 820.25
$ws.Cells.Item(102, 9).Value = 723.1429000000001
$ws.Cells.Item(102, 10).Value = 1500
$ws.Cells.Item(102, 11).Value = 723.1429000000001
$ws.Cells.Item(102, 12).Value = 1500
$ws.Cells.Item(102, 13).Value = 898.8570999999999
$ws.Cells.Item(102, 14).Value = -4744

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1897.1765
$ws.Cells.Item(100, 9).Value = 1778.8572
$ws.Cells.Item(100, 10).Value = 1980
$ws.Cells.Item(100, 11).Value = 1778.8572
$ws.Cells.Item(100, 12).Value = 1980
$ws.Cells.Item(100, 13).Value = -1237.8572
$ws.Cells.Item(100, 14).Value = -3062

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 61017.176
$ws.Cells.Item(132, 9).Value = 2406.3076
$ws.Cells.Item(132, 10).Value = 251502.5
$ws.Cells.Item(132, 11).Value = 7218.9228
$ws.Cells.Item(132, 12).Value = 754507.5
$ws.Cells.Item(132, 13).Value = -4688.9228
$ws.Cells.Item(132, 14).Value = -759567.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 96111.55
$ws.Cells.Item(136, 9).Value = 53392.316
$ws.Cells.Item(136, 11).Value = 160176.948
$ws.Cells.Item(136, 13).Value = -157626.948

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2209.1667
$ws.Cells.Item(122, 9).Value = 1585.6875
$ws.Cells.Item(122, 10).Value = 2921.7144
$ws.Cells.Item(122, 11).Value = 4757.0625
$ws.Cells.Item(122, 12).Value = 8765.143199999999
$ws.Cells.Item(122, 13).Value = -2307.0625
$ws.Cells.Item(122, 14).Value = -13665.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2119.6
$ws.Cells.Item(126, 9).Value = 2866
$ws.Cells.Item(126, 10).Value = 1000
$ws.Cells.Item(126, 11).Value = 8598
$ws.Cells.Item(126, 12).Value = 3000
$ws.Cells.Item(126, 13).Value = -6128
$ws.Cells.Item(126, 14).Value = -7940

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 43484.062
$ws.Cells.Item(132, 9).Value = 24653.596
$ws.Cells.Item(132, 10).Value = 201660
$ws.Cells.Item(132, 11).Value = 73960.788
$ws.Cells.Item(132, 12).Value = 604980
$ws.Cells.Item(132, 13).Value = -71430.788
$ws.Cells.Item(132, 14).Value = -610040

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 40892.844
$ws.Cells.Item(136, 9).Value = 24443.572
$ws.Cells.Item(136, 11).Value = 73330.716
$ws.Cells.Item(136, 13).Value = -70780.716
